$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix accent on "Cordoba" -> "Córdoba" (row 23, airport name column)
$ws.Range("B23").Value2 = "Córdoba"

# The footer notes at the bottom of the sheet are rearranged:
#  - old B85 "Fuente: AFAC. Agencia Federal de Aviación Civil." moves down to B86
#  - B85 now holds the updated "Actualización: mayo 2024." note
#  - old D86 "Ultima actualización: mayo 2024" is removed (cell deleted)
#  - old D87 "Dirección General de Planeación" text is removed, but the cell/style remains (now empty)

$ws.Range("B85").Value2 = "Actualización: mayo 2024."
$ws.Range("B86").Value2 = "Fuente: AFAC. Agencia Federal de Aviación Civil."

# Remove the D86 cell entirely (content + formatting)
$ws.Range("D86").Clear()

# Clear D87's text but keep its formatting/style
$ws.Range("D87").ClearContents()
